$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "connect1_1"
$ws.Range("B13").Value = 0.01
